$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 316, shifting existing rows 316:402 down to 317:403
$ws.Rows.Item(316).Insert()

# Populate the newly inserted row 316 with the new weekly data point
$ws.Range("A316").Value = 6
$ws.Range("B316").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C316").Value = "Metropolitana"
$ws.Range("D316").Value = 44642
$ws.Range("E316").Value = 13
$ws.Range("F316").Value = 100112039
$ws.Range("G316").Value = "Ciboulette"
$ws.Range("H316").Value = "Sin especificar"
$ws.Range("I316").Value = "Primera"
$ws.Range("J316").Value = 630
$ws.Range("K316").Value = 1300
$ws.Range("L316").Value = 1500
$ws.Range("M316").Value = 1379
$ws.Range("N316").Value = "$/docena de atados"
$ws.Range("O316").Value = "Región Metropolitana"
$ws.Range("P316").Value = 460
$ws.Range("Q316").Value = 3
$ws.Range("R316").Value = "Hortaliza"
